# The roster data for Darius Garland (row 7) and Donovan Mitchell (row 8)
# had been swapped by mistake; this corrects it by swapping the full row
# contents (every column except the sequential "No." column A) back into
# their proper places.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 7
$row2 = 8

$rng1 = $ws.Range("B$row1`:K$row1")
$rng2 = $ws.Range("B$row2`:K$row2")

# Use a scratch area far away from the live data and Copy/Paste (rather than
# reading/writing .Value) so that cell types (e.g. text values that look like
# numbers) are preserved exactly as they were, instead of being re-inferred.
$scratch = $ws.Range("B1000:K1000")

$rng1.Copy($scratch)
$rng2.Copy($rng1)
$scratch.Copy($rng2)
$scratch.Clear()
